$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 6 (shifts existing rows 6-52 down to 7-53)
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new weekly record
$ws.Range("A6").Value = 11
$ws.Range("B6").Value = "Vega Monumental Concepción"
$ws.Range("C6").Value = "Bíobío"
$ws.Range("D6").Value = 44552
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 100112012
$ws.Range("G6").Value = "Espinaca"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 60
$ws.Range("K6").Value = 11000
$ws.Range("L6").Value = 12000
$ws.Range("M6").Value = 11500
$ws.Range("N6").Value = "$/cuna 10 kilos"
$ws.Range("O6").Value = "Región Metropolitana"
$ws.Range("P6").Value = 1150
$ws.Range("Q6").Value = 10
$ws.Range("R6").Value = "Hortaliza"
